$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh adds two new price records (market date 2021-09-09,
# serial 44448) at the top of the existing "Pepino ensalada" block, pushing
# the previously-existing rows 249-254 down to 251-256.
$ws.Rows("249:250").Insert()

# New row 249: Primera quality
$ws.Range("A249").Value2 = 6
$ws.Range("B249").Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C249").Value2 = 'Metropolitana'
$ws.Range("D249").Value2 = 44448
$ws.Range("E249").Value2 = 13
$ws.Range("F249").Value2 = 100112043
$ws.Range("G249").Value2 = 'Pepino ensalada'
$ws.Range("H249").Value2 = 'Sin especificar'
$ws.Range("I249").Value2 = 'Primera'
$ws.Range("J249").Value2 = 400
$ws.Range("K249").Value2 = 16000
$ws.Range("L249").Value2 = 17000
$ws.Range("M249").Value2 = 16425
$ws.Range("N249").Value2 = '$/caja 60 unidades'
$ws.Range("O249").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P249").Value2 = 274
$ws.Range("Q249").Value2 = 60
$ws.Range("R249").Value2 = 'Hortaliza'

# New row 250: Segunda quality
$ws.Range("A250").Value2 = 6
$ws.Range("B250").Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C250").Value2 = 'Metropolitana'
$ws.Range("D250").Value2 = 44448
$ws.Range("E250").Value2 = 13
$ws.Range("F250").Value2 = 100112043
$ws.Range("G250").Value2 = 'Pepino ensalada'
$ws.Range("H250").Value2 = 'Sin especificar'
$ws.Range("I250").Value2 = 'Segunda'
$ws.Range("J250").Value2 = 200
$ws.Range("K250").Value2 = 13000
$ws.Range("L250").Value2 = 14000
$ws.Range("M250").Value2 = 13600
$ws.Range("N250").Value2 = '$/caja 100 unidades'
$ws.Range("O250").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P250").Value2 = 136
$ws.Range("Q250").Value2 = 100
$ws.Range("R250").Value2 = 'Hortaliza'
